$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the date-block (shifts existing rows 1119-1251 down to 1121-1253)
$ws.Rows("1119:1120").Insert()

# New row 1119: Primera quality, week of 2023-10-13 (serial 45212)
$ws.Range("A1119").Value = 3
$ws.Range("B1119").Value = "Femacal de La Calera"
$ws.Range("C1119").Value = "Coquimbo"
$ws.Range("D1119").Value = 45212
$ws.Range("E1119").Value = 5
$ws.Range("F1119").Value = 100114014
$ws.Range("G1119").Value = "Betarraga"
$ws.Range("H1119").Value = "Sin especificar"
$ws.Range("I1119").Value = "Primera"
$ws.Range("J1119").Value = 2800
$ws.Range("K1119").Value = 500
$ws.Range("L1119").Value = 550
$ws.Range("M1119").Value = 521
$ws.Range("N1119").Value = "$/paquete 4 unidades"
$ws.Range("O1119").Value = "Provincia de Quillota"
$ws.Range("P1119").Value = 130
$ws.Range("Q1119").Value = 4
$ws.Range("R1119").Value = "Hortaliza"

# New row 1120: Segunda quality, week of 2023-10-13 (serial 45212)
$ws.Range("A1120").Value = 3
$ws.Range("B1120").Value = "Femacal de La Calera"
$ws.Range("C1120").Value = "Coquimbo"
$ws.Range("D1120").Value = 45212
$ws.Range("E1120").Value = 5
$ws.Range("F1120").Value = 100114014
$ws.Range("G1120").Value = "Betarraga"
$ws.Range("H1120").Value = "Sin especificar"
$ws.Range("I1120").Value = "Segunda"
$ws.Range("J1120").Value = 1200
$ws.Range("K1120").Value = 400
$ws.Range("L1120").Value = 400
$ws.Range("M1120").Value = 400
$ws.Range("N1120").Value = "$/paquete 4 unidades"
$ws.Range("O1120").Value = "Provincia de Quillota"
$ws.Range("P1120").Value = 100
$ws.Range("Q1120").Value = 4
$ws.Range("R1120").Value = "Hortaliza"
